$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44902
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = "$/caja 10 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1550
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44902
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 13000
$ws.Range("Q3").Value = "$/caja 10 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1300
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44918
$ws.Range("K4").Value = "Dina"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 972
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44559
$ws.Range("K5").Value = "Modesto"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1083
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44559
$ws.Range("K6").Value = "Modesto"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = "$/caja 18 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44938
$ws.Range("K7").Value = "Modesto"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 270
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14556
$ws.Range("Q7").Value = "$/caja 15 kilos"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 970
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 44187
$ws.Range("K8").Value = "Dina"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15500
$ws.Range("Q8").Value = "$/caja 18 kilos"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 861
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44908
$ws.Range("K9").Value = "Albaricoque"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21000
$ws.Range("Q9").Value = "$/caja 18 kilos granel"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1167
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44189
$ws.Range("K10").Value = "Dina"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = "$/caja 15 kilos granel"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1033
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("D11").Value = 44189
$ws.Range("K11").Value = "Dina"
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("Q11").Value = "$/caja 15 kilos granel"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 933
$ws.Range("T11").Value = 15

# Row 12
$ws.Range("D12").Value = 44944
$ws.Range("K12").Value = "Modesto"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 16500
$ws.Range("Q12").Value = "$/caja 16 kilos empedrada"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 1031
$ws.Range("T12").Value = 16

# Row 13
$ws.Range("D13").Value = 44944
$ws.Range("K13").Value = "Modesto"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("Q13").Value = "$/caja 16 kilos empedrada"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 875
$ws.Range("T13").Value = 16

# Row 14
$ws.Range("D14").Value = 44875
$ws.Range("K14").Value = "Castle Brite"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 31000
$ws.Range("O14").Value = 32000
$ws.Range("P14").Value = 31400
$ws.Range("Q14").Value = "$/bandeja 10 kilos"
$ws.Range("R14").Value = "Provincia de Limarí"
$ws.Range("S14").Value = 3140
$ws.Range("T14").Value = 10

# Row 15
$ws.Range("D15").Value = 44545
$ws.Range("K15").Value = "Castle Brite"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 19000
$ws.Range("P15").Value = 18500
$ws.Range("Q15").Value = "$/caja 15 kilos"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 1233
$ws.Range("T15").Value = 15

# Row 16
$ws.Range("D16").Value = 44545
$ws.Range("K16").Value = "Castle Brite"
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 17000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 17000
$ws.Range("Q16").Value = "$/caja 15 kilos"
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 1133
$ws.Range("T16").Value = 15

# Row 17
$ws.Range("D17").Value = 44159
$ws.Range("K17").Value = "Castle Brite"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 14000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 14500
$ws.Range("Q17").Value = "$/caja 15 kilos"
$ws.Range("R17").Value = "Región Metropolitana"
$ws.Range("S17").Value = 967
$ws.Range("T17").Value = 15

# Row 18
$ws.Range("D18").Value = 44579
$ws.Range("K18").Value = "Modesto"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 180
$ws.Range("N18").Value = 13000
$ws.Range("O18").Value = 14000
$ws.Range("P18").Value = 13444
$ws.Range("Q18").Value = "$/caja 18 kilos"
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 747
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44901
$ws.Range("K19").Value = "Castle Brite"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 15000
$ws.Range("O19").Value = 16000
$ws.Range("P19").Value = 15500
$ws.Range("Q19").Value = "$/caja 10 kilos"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 1550
$ws.Range("T19").Value = 10
